# Applies numeric updates to LevePriceNQ/HQ-related columns (H-N) across several
# sheets, as scraped from the scheduled price-refresh run.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 17 (Leve Item ID 38956)
$ws.Range("H17").Value = 3489
$ws.Range("J17").Value = 3489
$ws.Range("L17").Value = 10467
$ws.Range("N17").Value = -10803

# Row 19 (Leve Item ID 7015)
$ws.Range("H19").Value = 1484.6
$ws.Range("I19").Value = 1627.1
$ws.Range("J19").Value = 1199.6
$ws.Range("K19").Value = 1627.1
$ws.Range("L19").Value = 1199.6
$ws.Range("M19").Value = -1452.1
$ws.Range("N19").Value = -1549.6

# Row 28 (Leve Item ID 27772)
$ws.Range("H28").Value = 347.2
$ws.Range("I28").Value = 347.2
$ws.Range("K28").Value = 347.2
$ws.Range("M28").Value = 137.8

# Row 62 (Leve Item ID 27781)
$ws.Range("H62").Value = 5000
$ws.Range("J62").Value = 5000
$ws.Range("L62").Value = 5000
$ws.Range("N62").Value = -6248

# Row 65 (Leve Item ID 27781)
$ws.Range("H65").Value = 5000
$ws.Range("J65").Value = 5000
$ws.Range("L65").Value = 25000
$ws.Range("N65").Value = -31240

# Row 100 (Leve Item ID 19906)
$ws.Range("H100").Value = 625
$ws.Range("I100").Value = 750
$ws.Range("K100").Value = 750
$ws.Range("M100").Value = -209

# Row 107 (Leve Item ID 27766)
$ws.Range("H107").Value = 183.42857
$ws.Range("I107").Value = 172.44444
$ws.Range("J107").Value = 203.2
$ws.Range("K107").Value = 172.44444
$ws.Range("L107").Value = 203.2
$ws.Range("M107").Value = 1747.55556
$ws.Range("N107").Value = -4043.2

# Row 132 (Leve Item ID 44049)
$ws.Range("H132").Value = 883.0833
$ws.Range("I132").Value = 509.72726
$ws.Range("K132").Value = 1529.18178
$ws.Range("M132").Value = 1000.81822

# Row 137 (Leve Item ID 44013)
$ws.Range("H137").Value = 1050.7273
$ws.Range("I137").Value = 794.875
$ws.Range("J137").Value = 1733
$ws.Range("K137").Value = 2384.625
$ws.Range("L137").Value = 5199
$ws.Range("M137").Value = 165.375
$ws.Range("N137").Value = -10299

$ws = $wb.Worksheets.Item("ARM")
# Row 2 (Leve Item ID 27713)
$ws.Range("H2").Value = 1133.3334
$ws.Range("I2").Value = 1133.3334
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 1133.3334
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = -1020.3334
$ws.Range("N2").ClearContents()

# Row 5 (Leve Item ID 5091)
$ws.Range("H5").Value = 608.3333
$ws.Range("J5").Value = 233.33333
$ws.Range("L5").Value = 233.33333
$ws.Range("N5").Value = -457.33333

# Row 32 (Leve Item ID 44147)
$ws.Range("H32").Value = 5094.5
$ws.Range("I32").Value = 5094.5
$ws.Range("K32").Value = 5094.5
$ws.Range("M32").Value = -4807.5

# Row 61 (Leve Item ID 43999)
$ws.Range("H61").Value = 2553.5881
$ws.Range("I61").Value = 1321.2
$ws.Range("J61").Value = 4314.143
$ws.Range("K61").Value = 1321.2
$ws.Range("L61").Value = 4314.143
$ws.Range("M61").Value = -1109.2
$ws.Range("N61").Value = -4738.143

# Row 110 (Leve Item ID 27708)
$ws.Range("H110").Value = 1200.7778
$ws.Range("I110").Value = 1215.4286
$ws.Range("J110").Value = 1149.5
$ws.Range("K110").Value = 1215.4286
$ws.Range("L110").Value = 1149.5
$ws.Range("M110").Value = 829.5714
$ws.Range("N110").Value = -5239.5

# Row 116 (Leve Item ID 27713)
$ws.Range("H116").Value = 1133.3334
$ws.Range("I116").Value = 1133.3334
$ws.Range("J116").Value = 0
$ws.Range("K116").Value = 1133.3334
$ws.Range("L116").Value = 0
$ws.Range("M116").Value = 1160.6666
$ws.Range("N116").ClearContents()

# Row 122 (Leve Item ID 36168)
$ws.Range("H122").Value = 1864.9412
$ws.Range("I122").Value = 1191.7273
$ws.Range("K122").Value = 3575.1819
$ws.Range("M122").Value = -1125.1819

# Row 132 (Leve Item ID 43997)
$ws.Range("H132").Value = 3061.4666
$ws.Range("I132").Value = 1365.5
$ws.Range("K132").Value = 4096.5
$ws.Range("M132").Value = -1566.5

# Row 136 (Leve Item ID 43999)
$ws.Range("H136").Value = 2553.5881
$ws.Range("I136").Value = 1321.2
$ws.Range("J136").Value = 4314.143
$ws.Range("K136").Value = 3963.6
$ws.Range("L136").Value = 12942.429
$ws.Range("M136").Value = -1413.6
$ws.Range("N136").Value = -18042.429

$ws = $wb.Worksheets.Item("BSM")
# Row 3 (Leve Item ID 27713)
$ws.Range("H3").Value = 1133.3334
$ws.Range("I3").Value = 1133.3334
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 1133.3334
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = -1019.3334
$ws.Range("N3").ClearContents()

# Row 4 (Leve Item ID 5091)
$ws.Range("H4").Value = 608.3333
$ws.Range("J4").Value = 233.33333
$ws.Range("L4").Value = 233.33333
$ws.Range("N4").Value = -463.33333

$ws = $wb.Worksheets.Item("CRP")
# Row 31 (Leve Item ID 44023)
$ws.Range("H31").Value = 2669.1482
$ws.Range("I31").Value = 1604.375
$ws.Range("K31").Value = 1604.375
$ws.Range("M31").Value = -1309.375

# Row 34 (Leve Item ID 44023)
$ws.Range("H34").Value = 2669.1482
$ws.Range("I34").Value = 1604.375
$ws.Range("K34").Value = 1604.375
$ws.Range("M34").Value = -1402.375

# Row 99 (Leve Item ID 36198)
$ws.Range("H99").Value = 3918.2
$ws.Range("I99").Value = 3918.2
$ws.Range("K99").Value = 3918.2
$ws.Range("M99").Value = -2420.2

# Row 107 (Leve Item ID 27689)
$ws.Range("H107").Value = 837.8
$ws.Range("I107").Value = 863.3333
$ws.Range("J107").Value = 799.5
$ws.Range("K107").Value = 863.3333
$ws.Range("L107").Value = 799.5
$ws.Range("M107").Value = 1056.6667
$ws.Range("N107").Value = -4639.5

# Row 126 (Leve Item ID 36198)
$ws.Range("H126").Value = 3918.2
$ws.Range("I126").Value = 3918.2
$ws.Range("K126").Value = 11754.6
$ws.Range("M126").Value = -9284.599999999999

# Row 132 (Leve Item ID 44019)
$ws.Range("H132").Value = 2796.2307
$ws.Range("I132").Value = 2261.5557
$ws.Range("J132").Value = 3999.25
$ws.Range("K132").Value = 6784.6671
$ws.Range("L132").Value = 11997.75
$ws.Range("M132").Value = -4254.6671
$ws.Range("N132").Value = -17057.75

# Row 134 (Leve Item ID 44020)
$ws.Range("H134").Value = 2836.875
$ws.Range("I134").Value = 2666.8333
$ws.Range("J134").Value = 3347
$ws.Range("K134").Value = 8000.499899999999
$ws.Range("L134").Value = 10041
$ws.Range("M134").Value = -5465.499899999999
$ws.Range("N134").Value = -15111

$ws = $wb.Worksheets.Item("CUL")
# Row 8 (Leve Item ID 16734)
$ws.Range("H8").Value = 918.5714
$ws.Range("I8").Value = 918.5714
$ws.Range("K8").Value = 2755.7142
$ws.Range("M8").Value = -2616.7142

# Row 108 (Leve Item ID 27853)
$ws.Range("H108").Value = 1989.75
$ws.Range("I108").Value = 1989.75
$ws.Range("K108").Value = 5969.25
$ws.Range("M108").Value = -3089.25

$ws = $wb.Worksheets.Item("GSM")
# Row 2 (Leve Item ID 5062)
$ws.Range("H2").Value = 154.5
$ws.Range("I2").Value = 148
$ws.Range("K2").Value = 148
$ws.Range("M2").Value = -35

# Row 132 (Leve Item ID 44008)
$ws.Range("H132").Value = 2732.5833
$ws.Range("I132").Value = 2224.125
$ws.Range("J132").Value = 3749.5
$ws.Range("K132").Value = 6672.375
$ws.Range("L132").Value = 11248.5
$ws.Range("M132").Value = -4142.375
$ws.Range("N132").Value = -16308.5

$ws = $wb.Worksheets.Item("LTW")
# Row 40 (Leve Item ID 36248)
$ws.Range("H40").Value = 833.3333
$ws.Range("I40").Value = 833.3333
$ws.Range("K40").Value = 833.3333
$ws.Range("M40").Value = -697.3333

# Row 132 (Leve Item ID 44058)
$ws.Range("H132").Value = 2317.5264
$ws.Range("I132").Value = 1982.3334
$ws.Range("J132").Value = 3574.5
$ws.Range("K132").Value = 5947.0002
$ws.Range("L132").Value = 10723.5
$ws.Range("M132").Value = -3417.0002
$ws.Range("N132").Value = -15783.5

# Row 136 (Leve Item ID 44060)
$ws.Range("H136").Value = 2542.8948
$ws.Range("I136").Value = 2573.111
$ws.Range("K136").Value = 7719.333
$ws.Range("M136").Value = -5169.333

$ws = $wb.Worksheets.Item("WVR")
# Row 122 (Leve Item ID 36208)
$ws.Range("H122").Value = 3001
$ws.Range("I122").Value = 2999
$ws.Range("J122").Value = 3001.6667
$ws.Range("K122").Value = 8997
$ws.Range("L122").Value = 9005.000100000001
$ws.Range("M122").Value = -6547
$ws.Range("N122").Value = -13905.0001

# Row 132 (Leve Item ID 44029)
$ws.Range("H132").Value = 2000.1936
$ws.Range("I132").Value = 846.2222
$ws.Range("J132").Value = 3598
$ws.Range("K132").Value = 2538.6666
$ws.Range("L132").Value = 10794
$ws.Range("M132").Value = -8.666600000000017
$ws.Range("N132").Value = -15854

# Row 136 (Leve Item ID 44031)
$ws.Range("H136").Value = 899.6667
$ws.Range("I136").Value = 899.6667
$ws.Range("K136").Value = 2699.0001
$ws.Range("M136").Value = -149.0001000000002
